$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 8.979873666666666
$ws.Range("N2").Value = 26.939621
$ws.Range("O2").Value = 0.3651616045144693
$ws.Range("P2").Value = 0.3651616045144694
$ws.Range("Q2").Value = 2.021342622745666
$ws.Range("R2").Value = 18.192083604711
$ws.Range("S2").Value = 0.07011789852871207
$ws.Range("T2").Value = 0.07011789852871209

# Row 3
$ws.Range("M3").Value = 4.482719
$ws.Range("O3").Value = 0.1822872930499836
$ws.Range("P3").Value = 0.1822872930499837
$ws.Range("S3").Value = 0.03500258997423123
$ws.Range("T3").Value = 0.03500258997423124

# Row 4
$ws.Range("M4").Value = 9.285498666666667
$ws.Range("N4").Value = 27.856496
$ws.Range("O4").Value = 0.3775896763919173
$ws.Range("P4").Value = 0.3775896763919173
$ws.Range("Q4").Value = 2.090137893370667
$ws.Range("R4").Value = 18.811241040336
$ws.Range("S4").Value = 0.07250432216152834
$ws.Range("T4").Value = 0.07250432216152834

# Row 5
$ws.Range("M5").Value = 1.843414333333333
$ws.Range("N5").Value = 5.530243
$ws.Range("O5").Value = 0.07496142604362967
$ws.Range("P5").Value = 0.07496142604362969
$ws.Range("Q5").Value = 0.4149470361903334
$ws.Range("R5").Value = 3.734523325713
$ws.Range("S5").Value = 0.01439400418859346
$ws.Range("T5").Value = 0.01439400418859346

# Row 6
$ws.Range("G6").Value = 0.9471683333333334
$ws.Range("I6").Value = 0.807981185146935
$ws.Range("M6").Value = 8.979873666666666
$ws.Range("N6").Value = 26.939621
$ws.Range("O6").Value = 0.3651616045144693
$ws.Range("P6").Value = 0.3651616045144694
$ws.Range("Q6").Value = 8.505451974400556
$ws.Range("R6").Value = 76.549067769605
$ws.Range("S6").Value = 0.2950437059857573
$ws.Range("T6").Value = 0.2950437059857574

# Row 7
$ws.Range("G7").Value = 0.9471683333333334
$ws.Range("I7").Value = 0.807981185146935
$ws.Range("M7").Value = 4.482719
$ws.Range("O7").Value = 0.1822872930499836
$ws.Range("P7").Value = 0.1822872930499837
$ws.Range("Q7").Value = 4.245889484031667
$ws.Range("S7").Value = 0.1472847030757524
$ws.Range("T7").Value = 0.1472847030757524

# Row 8
$ws.Range("G8").Value = 0.9471683333333334
$ws.Range("I8").Value = 0.807981185146935
$ws.Range("M8").Value = 9.285498666666667
$ws.Range("N8").Value = 27.856496
$ws.Range("O8").Value = 0.3775896763919173
$ws.Range("P8").Value = 0.3775896763919173
$ws.Range("Q8").Value = 8.794930296275556
$ws.Range("R8").Value = 79.15437266648
$ws.Range("S8").Value = 0.305085354230389
$ws.Range("T8").Value = 0.305085354230389

# Row 9
$ws.Range("G9").Value = 0.9471683333333334
$ws.Range("I9").Value = 0.807981185146935
$ws.Range("M9").Value = 1.843414333333333
$ws.Range("N9").Value = 5.530243
$ws.Range("O9").Value = 0.07496142604362967
$ws.Range("P9").Value = 0.07496142604362969
$ws.Range("Q9").Value = 1.746023681746111
$ws.Range("R9").Value = 15.714213135715
$ws.Range("S9").Value = 0.06056742185503621
$ws.Range("T9").Value = 0.06056742185503623
